# combined Add_Lottery_Admin and Add_Lottery_Report into one script
# Adds three new user rows (15, 16, 17) to the "Users" sheet, each with a
# mailto hyperlink on the Email column, matching the existing rows' layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Row 15: Sam Amanyu - Lottery Admin
$ws.Range("A15").Value = "samanyu5@jolongestr.com"
$ws.Range("B15").Value = "Sam"
$ws.Range("C15").Value = "Amanyu"
$ws.Range("D15").Value = 6612200748
$ws.Range("E15").Value = "Lottery Admin"
$ws.Hyperlinks.Add($ws.Range("A15"), "mailto:samanyu5@jolongestr.com")
$ws.Range("A15").Style = "Hyperlink"

# Row 16: Mary Alia - Lottery Report
$ws.Range("A16").Value = "malia25@jolongestr.com"
$ws.Range("B16").Value = "Mary"
$ws.Range("C16").Value = "Alia"
$ws.Range("D16").Value = 6612200748
$ws.Range("E16").Value = "Lottery Report"
$ws.Hyperlinks.Add($ws.Range("A16"), "mailto:malia25@jolongestr.com")
$ws.Range("A16").Style = "Hyperlink"

# Row 17: Sam Amanyu - Lottery Admin (second account)
$ws.Range("A17").Value = "samantha5@jolongestr.com"
$ws.Range("B17").Value = "Sam"
$ws.Range("C17").Value = "Amanyu"
$ws.Range("D17").Value = 6612200748
$ws.Range("E17").Value = "Lottery Admin"
$ws.Hyperlinks.Add($ws.Range("A17"), "mailto:samantha5@jolongestr.com")
$ws.Range("A17").Style = "Hyperlink"

# Move selection to match the author's final cursor position
$ws.Range("E19").Select()
